# Fruta / hortaliza, semanal
#
# Two new weekly price records (Especial / Primera, Kiwi Hayward,
# Macroferia Regional de Talca) were captured for 2022-08-09 and need to
# be inserted into the dataset at row 246, ahead of the existing history
# (which is sorted most-recent-first). Inserting the two rows pushes the
# remainder of the table (old rows 246-311) down by two rows, which is
# exactly what the workbook's new A1:T313 dimension reflects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 246-311 down two rows, opening up space for the
# two new records at rows 246:247.
$ws.Rows("246:247").Insert()

# New records for 2022-08-09.
$newRows = New-Object 'object[,]' 2,20

$newRows[0,0]  = 5
$newRows[0,1]  = "Macroferia Regional de Talca"
$newRows[0,2]  = "Maule"
$newRows[0,3]  = (Get-Date -Year 2022 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0)
$newRows[0,4]  = 7
$newRows[0,5]  = "Fruta"
$newRows[0,6]  = 100101
$newRows[0,7]  = "Berries"
$newRows[0,8]  = 100101007
$newRows[0,9]  = "Kiwi"
$newRows[0,10] = "Hayward"
$newRows[0,11] = "Especial"
$newRows[0,12] = 200
$newRows[0,13] = 8000
$newRows[0,14] = 8000
$newRows[0,15] = 8000
$newRows[0,16] = "`$/bandeja 18 kilos"
$newRows[0,17] = "Provincia de Curicó"
$newRows[0,18] = 444
$newRows[0,19] = 18

$newRows[1,0]  = 5
$newRows[1,1]  = "Macroferia Regional de Talca"
$newRows[1,2]  = "Maule"
$newRows[1,3]  = (Get-Date -Year 2022 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0)
$newRows[1,4]  = 7
$newRows[1,5]  = "Fruta"
$newRows[1,6]  = 100101
$newRows[1,7]  = "Berries"
$newRows[1,8]  = 100101007
$newRows[1,9]  = "Kiwi"
$newRows[1,10] = "Hayward"
$newRows[1,11] = "Primera"
$newRows[1,12] = 190
$newRows[1,13] = 6000
$newRows[1,14] = 6000
$newRows[1,15] = 6000
$newRows[1,16] = "`$/bandeja 18 kilos"
$newRows[1,17] = "Provincia de Curicó"
$newRows[1,18] = 333
$newRows[1,19] = 18

$ws.Range("A246:T247").Value = $newRows

# Match the date-column style (s="2", "YYYY-MM-DD HH:MM:SS") used by the
# rest of the Fecha column.
$ws.Range("D246:D247").NumberFormat = $ws.Range("D248").NumberFormat

Write-Output "Inserted 2 new rows (246:247) and shifted existing data down; dimension now A1:T313."
